$d = $word.ActiveDocument

# Locate the list-item paragraphs whose entire text is "City_name" or
# "State_name" (the sub-bullets under "Remove these columns") and remove
# them completely, including their paragraph marks, so the list collapses
# from Area_name / City_name / State_name / Country_name down to just
# Area_name / Country_name.
$ranges = @()
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.Trim()
    if ($t -eq "City_name" -or $t -eq "State_name") {
        $ranges += ,@($p.Range.Start, $p.Range.End)
    }
}

# Delete from the end of the document backwards so earlier offsets stay
# valid while later ones are removed.
$sorted = $ranges | Sort-Object -Property { $_[0] } -Descending
foreach ($r in $sorted) {
    $rng = $d.Range($r[0], $r[1])
    $rng.Delete()
}
